# Updated symbol list data (Price / Volume(1h)) per commit on 2023-01-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "303.57"
$cell.ClearFormats()
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "4.84%"
$cell.ClearFormats()
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "35.83"
$cell.ClearFormats()
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "15.33%"
$cell.ClearFormats()
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "5.151"
$cell.ClearFormats()
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "4.07%"
$cell.ClearFormats()
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "0.07867"
$cell.ClearFormats()
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "6.94%"
$cell.ClearFormats()
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "-1.04%"
$cell.ClearFormats()
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "8.056"
$cell.ClearFormats()
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "4.92%"
$cell.ClearFormats()
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "3.978"
$cell.ClearFormats()
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "6.26%"
$cell.ClearFormats()
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.9280"
$cell.ClearFormats()
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "1.02%"
$cell.ClearFormats()
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.1007"
$cell.ClearFormats()
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "10.53%"
$cell.ClearFormats()
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "7.29%"
$cell.ClearFormats()
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.08535"
$cell.ClearFormats()
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "4.39%"
$cell.ClearFormats()
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "8.56%"
$cell.ClearFormats()
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.09932"
$cell.ClearFormats()
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "-0.50%"
$cell.ClearFormats()
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.001494"
$cell.ClearFormats()
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "-0.04%"
$cell.ClearFormats()
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.005741"
$cell.ClearFormats()
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "0.36%"
$cell.ClearFormats()
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "3.484"
$cell.ClearFormats()
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "4.81%"
$cell.ClearFormats()
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "3.05%"
$cell.ClearFormats()
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.1326"
$cell.ClearFormats()
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "2.04%"
$cell.ClearFormats()
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "4.520"
$cell.ClearFormats()
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "7.91%"
$cell.ClearFormats()
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.2218"
$cell.ClearFormats()
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "4.39%"
$cell.ClearFormats()
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.04635"
$cell.ClearFormats()
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "2.64%"
$cell.ClearFormats()
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "0.001217"
$cell.ClearFormats()
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "0.08%"
$cell.ClearFormats()
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.004494"
$cell.ClearFormats()
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "7.08%"
$cell.ClearFormats()
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.0001295"
$cell.ClearFormats()
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "-0.41%"
$cell.ClearFormats()
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.0003392"
$cell.ClearFormats()
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "-0.09%"
$cell.ClearFormats()
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.01745"
$cell.ClearFormats()
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "10.35%"
$cell.ClearFormats()
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.04725"
$cell.ClearFormats()
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "4.92%"
$cell.ClearFormats()
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.007813"
$cell.ClearFormats()
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "5.92%"
$cell.ClearFormats()
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.1421"
$cell.ClearFormats()
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "6.44%"
$cell.ClearFormats()
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.008802"
$cell.ClearFormats()
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "-10.75%"
$cell.ClearFormats()
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.002211"
$cell.ClearFormats()
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "-0.41%"
$cell.ClearFormats()
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.009150"
$cell.ClearFormats()
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "7.79%"
$cell.ClearFormats()
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.00006062"
$cell.ClearFormats()
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "-0.87%"
$cell.ClearFormats()
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.00000000749"
$cell.ClearFormats()
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "-0.18%"
$cell.ClearFormats()
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "5.799"
$cell.ClearFormats()
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "137.15%"
$cell.ClearFormats()
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.002685"
$cell.ClearFormats()
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "34.23%"
$cell.ClearFormats()
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.00002096"
$cell.ClearFormats()
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "-0.18%"
$cell.ClearFormats()
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.0001996"
$cell.ClearFormats()
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "-0.18%"
$cell.ClearFormats()
